# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the "Leve profit" sheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5169102
$ws.Range("J17").Value = 6137602
$ws.Range("L17").Value = 18412806
$ws.Range("N17").Value = -18413142
$ws.Range("H21").Value = 10999.5
$ws.Range("I21").Value = 10999.5
$ws.Range("K21").Value = 10999.5
$ws.Range("M21").Value = -10531.5
$ws.Range("H23").Value = 10999.5
$ws.Range("I23").Value = 10999.5
$ws.Range("K23").Value = 10999.5
$ws.Range("M23").Value = -10765.5
$ws.Range("H40").Value = 1711.1111
$ws.Range("I40").Value = 1673.3334
$ws.Range("K40").Value = 1673.3334
$ws.Range("M40").Value = -1498.3334
$ws.Range("I42").Value = 68.125
$ws.Range("J42").Value = 31.333334
$ws.Range("K42").Value = 204.375
$ws.Range("L42").Value = 94.00000199999999
$ws.Range("M42").Value = 25.625
$ws.Range("N42").Value = -554.000002
$ws.Range("H51").Value = 7736.875
$ws.Range("J51").Value = 9974.25
$ws.Range("L51").Value = 9974.25
$ws.Range("N51").Value = -10942.25
$ws.Range("H86").Value = 6984.3335
$ws.Range("I86").Value = 3967.3333
$ws.Range("J86").Value = 10001.333
$ws.Range("K86").Value = 3967.3333
$ws.Range("L86").Value = 10001.333
$ws.Range("M86").Value = -2844.3333
$ws.Range("N86").Value = -12247.333
$ws.Range("H89").Value = 6984.3335
$ws.Range("I89").Value = 3967.3333
$ws.Range("J89").Value = 10001.333
$ws.Range("K89").Value = 19836.6665
$ws.Range("L89").Value = 50006.665
$ws.Range("M89").Value = -14220.6665
$ws.Range("N89").Value = -61238.665
$ws.Range("H100").Value = 44586.4
$ws.Range("I100").Value = 75837.36
$ws.Range("K100").Value = 75837.36
$ws.Range("M100").Value = -75296.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4137
$ws.Range("I32").Value = 2140.258
$ws.Range("J32").Value = 11874.375
$ws.Range("K32").Value = 2140.258
$ws.Range("L32").Value = 11874.375
$ws.Range("M32").Value = -1853.258
$ws.Range("N32").Value = -12448.375
$ws.Range("H97").Value = 1032.7
$ws.Range("I97").Value = 1080.8334
$ws.Range("K97").Value = 1080.8334
$ws.Range("M97").Value = -584.8334
$ws.Range("H122").Value = 1375.6666
$ws.Range("I122").Value = 1265.0476
$ws.Range("K122").Value = 3795.142800000001
$ws.Range("M122").Value = -1345.142800000001
$ws.Range("H132").Value = 3061.6365
$ws.Range("I132").Value = 3056.907
$ws.Range("K132").Value = 9170.721000000001
$ws.Range("M132").Value = -6640.721000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 922.2143
$ws.Range("I94").Value = 761.3
$ws.Range("K94").Value = 761.3
$ws.Range("M94").Value = -310.3
$ws.Range("H134").Value = 12935.8
$ws.Range("I134").Value = 5370.067
$ws.Range("J134").Value = 24284.4
$ws.Range("K134").Value = 16110.201
$ws.Range("L134").Value = 72853.20000000001
$ws.Range("M134").Value = -13575.201
$ws.Range("N134").Value = -77923.20000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2643.93
$ws.Range("J31").Value = 3330.1794
$ws.Range("L31").Value = 3330.1794
$ws.Range("N31").Value = -3920.1794
$ws.Range("H34").Value = 2643.93
$ws.Range("J34").Value = 3330.1794
$ws.Range("L34").Value = 3330.1794
$ws.Range("N34").Value = -3734.1794
$ws.Range("H41").Value = 12714.143
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H59").Value = 66273.84
$ws.Range("J59").Value = 67894.414
$ws.Range("L59").Value = 67894.414
$ws.Range("N59").Value = -70184.414
$ws.Range("H62").Value = 378673.75
$ws.Range("J62").Value = 503500.5
$ws.Range("L62").Value = 503500.5
$ws.Range("N62").Value = -504748.5
$ws.Range("H65").Value = 378673.75
$ws.Range("J65").Value = 503500.5
$ws.Range("L65").Value = 2517502.5
$ws.Range("N65").Value = -2523742.5
$ws.Range("H68").Value = 40000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 40000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 24000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 24000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H94").Value = 952.1579
$ws.Range("I94").Value = 481.2
$ws.Range("J94").Value = 1120.3572
$ws.Range("K94").Value = 481.2
$ws.Range("L94").Value = 1120.3572
$ws.Range("M94").Value = -30.19999999999999
$ws.Range("N94").Value = -2022.3572
$ws.Range("H103").Value = 32250
$ws.Range("I103").Value = 32250
$ws.Range("K103").Value = 32250
$ws.Range("M103").Value = -31078

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45880108
$ws.Range("I4").Value = 53526460
$ws.Range("K4").Value = 160579380
$ws.Range("M4").Value = -160579268
$ws.Range("H96").Value = 12500
$ws.Range("J96").Value = 12500
$ws.Range("L96").Value = 37500
$ws.Range("N96").Value = -41618
$ws.Range("H109").Value = 2999.2
$ws.Range("I109").Value = 2374
$ws.Range("J109").Value = 5500
$ws.Range("K109").Value = 7122
$ws.Range("L109").Value = 16500
$ws.Range("M109").Value = -6082
$ws.Range("N109").Value = -18580
$ws.Range("H117").Value = 928.1429000000001
$ws.Range("I117").Value = 624.5
$ws.Range("J117").Value = 1333
$ws.Range("K117").Value = 1873.5
$ws.Range("L117").Value = 3999
$ws.Range("M117").Value = 1568.5
$ws.Range("N117").Value = -10883
$ws.Range("H120").Value = 200
$ws.Range("I120").Value = 200
$ws.Range("K120").Value = 600
$ws.Range("M120").Value = 4238

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 2274.375
$ws.Range("I22").Value = 1956.4286
$ws.Range("K22").Value = 1956.4286
$ws.Range("M22").Value = -1427.4286
$ws.Range("H25").Value = 2166.3333
$ws.Range("J25").Value = 2749.5
$ws.Range("L25").Value = 2749.5
$ws.Range("N25").Value = -3807.5
$ws.Range("H102").Value = 2385.1177
$ws.Range("I102").Value = 2442.7144
$ws.Range("K102").Value = 2442.7144
$ws.Range("M102").Value = -820.7143999999998
$ws.Range("H113").Value = 335416.16
$ws.Range("I113").Value = 2499.75
$ws.Range("K113").Value = 2499.75
$ws.Range("M113").Value = -329.75
$ws.Range("H126").Value = 3518.4
$ws.Range("I126").Value = 3720.6667
$ws.Range("J126").Value = 1698
$ws.Range("K126").Value = 11162.0001
$ws.Range("L126").Value = 5094
$ws.Range("M126").Value = -8692.000100000001
$ws.Range("N126").Value = -10034

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1599.2858
$ws.Range("I46").Value = 1066.3334
$ws.Range("K46").Value = 1066.3334
$ws.Range("M46").Value = -878.3334
$ws.Range("H55").Value = 174.22223
$ws.Range("I55").Value = 92.36364
$ws.Range("J55").Value = 302.85715
$ws.Range("K55").Value = 92.36364
$ws.Range("L55").Value = 302.85715
$ws.Range("M55").Value = 80.63636
$ws.Range("N55").Value = -648.85715
$ws.Range("H63").Value = 45035.57
$ws.Range("I63").Value = 59999.668
$ws.Range("K63").Value = 59999.668
$ws.Range("M63").Value = -59250.668
$ws.Range("H66").Value = 45035.57
$ws.Range("I66").Value = 59999.668
$ws.Range("K66").Value = 179999.004
$ws.Range("M66").Value = -176255.004
$ws.Range("H68").Value = 2696.4082
$ws.Range("I68").Value = 2374.4524
$ws.Range("J68").Value = 4628.143
$ws.Range("K68").Value = 2374.4524
$ws.Range("L68").Value = 4628.143
$ws.Range("M68").Value = -1625.4524
$ws.Range("N68").Value = -6126.143
$ws.Range("H71").Value = 2696.4082
$ws.Range("I71").Value = 2374.4524
$ws.Range("J71").Value = 4628.143
$ws.Range("K71").Value = 11872.262
$ws.Range("L71").Value = 23140.715
$ws.Range("M71").Value = -8128.262000000001
$ws.Range("N71").Value = -30628.715
$ws.Range("H100").Value = 4642.625
$ws.Range("I100").Value = 3490.1667
$ws.Range("K100").Value = 3490.1667
$ws.Range("M100").Value = -2949.1667
$ws.Range("H132").Value = 3173.3635
$ws.Range("I132").Value = 2770.4243
$ws.Range("K132").Value = 8311.2729
$ws.Range("M132").Value = -5781.2729

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 28006
$ws.Range("J21").Value = 28006
$ws.Range("L21").Value = 28006
$ws.Range("N21").Value = -28476
$ws.Range("H23").Value = 230
$ws.Range("I23").Value = 62.5
$ws.Range("K23").Value = 62.5
$ws.Range("M23").Value = 166.5
$ws.Range("H24").Value = 6729.2
$ws.Range("J24").Value = 6729.2
$ws.Range("L24").Value = 6729.2
$ws.Range("N24").Value = -7189.2
$ws.Range("H34").Value = 1000
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H35").Value = 28006
$ws.Range("J35").Value = 28006
$ws.Range("L35").Value = 28006
$ws.Range("N35").Value = -28586
